$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A labels (status descriptions) to reflect the new breakdown
$ws.Range("A3").Value = "Cuenta con Placa Anterior y Adeudos"
$ws.Range("A4").Value = "Cuenta con Placa Anterior y Sin adeudos"
$ws.Range("A5").Value = "El número de placa no se localizó en el padrón"
$ws.Range("A6").Value = "EL VEHICULO FUE LOCALIZADO EN EL PADRON FISCAL DEL DISTRITO FEDERAL CON ESTATUS BAJA, ES NECESARIO QUE ACUDA A LA ADMINISTRACION TRIBUTARIA O CENTRO DE SERVICIO DE LA TESORERIA MAS CERCANO A SU DOMICILIO, A FIN DE REGISTRAR LOS DATOS DE SU VEHICULO, PRESENTANDO LOS SIGUIENTES DOCUMENTOS EN ORIGINAL Y COPIA PARA SU COTEJO"
$ws.Range("A7").Value = "PLACA CON PROBLEMAS DE ADEUDOS DEL IMPUESTO SOBRE TENENCIA, POR LO QUE ES NECESARIO QUE EL PROPIETARIO ACUDA EXCLUSIVAMENTE A LA ADMINISTRACIÓN AUXILIAR DE CENTRO HISTÓRICO, UBICADA EN IZAZAGA 89, MEZZANINE, COL. CENTRO, A FIN DE ACLARAR SU SITUACIÓN FISCAL, PRESENTANDO LOS RECIBOS DE PAGO ORIGINALES POR LOS EJERCICIOS FISCALES DE 2005 A 2012."
$ws.Range("A8").Value = "Sin adeudos"
$ws.Range("A9").Value = "VEHÍCULO CON ADEUDOS DE TENENCIA, FAVOR DE ACUDIR A LA ADMINISTRACIÓN TRIBUTARIA MÁS CERCANA A SU DOMICILIO, DE LUNES A VIERNES EN UN HORARIO DE 9:00 A 13:30 HORAS, CON LA SIGUIENTE DOCUMENTACIÓN EN ORIGINAL Y COPIA:"
$ws.Range("A10").Value = "Vehiculo reportado por Fiscalizacion y no puede verificar"

# Update column B counts
$ws.Range("B2").Value = 1150582
$ws.Range("B3").Value = 303941
$ws.Range("B4").Value = 155141
$ws.Range("B5").Value = 72182
$ws.Range("B6").Value = 77384
$ws.Range("B7").Value = 125
$ws.Range("B8").Value = 813613
$ws.Range("B9").Value = 1023
$ws.Range("B10").Value = 1118
